$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp column (O) for all data rows (2 through 398)
# from "2023-01-09 16:00:38" to "2023-01-09 20:49:46"
for ($r = 2; $r -le 398; $r++) {
    $cell = $ws.Cells.Item($r, 15)
    if ($cell.Value2 -eq "2023-01-09 16:00:38") {
        $cell.Value = "2023-01-09 20:49:46"
    }
}

# Update specific productAriaLabel (M) text to reflect "Online kein Bestand"
$ws.Cells.Item(247, 13).Value = "Pasquier Pitch Choco Barre Lait - Online kein Bestand 5.20 Schweizer Franken"
$ws.Cells.Item(313, 13).Value = "Mini Chococremecake - Online kein Bestand 4.20 Schweizer Franken"
